$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the header row with the new student-sample column names
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Roll No"
$ws.Range("C1").Value = "Father Name"
$ws.Range("D1").Value = "Course "
$ws.Range("E1").Value = "Session"
$ws.Range("F1").Value = "College"

# Remove the old trailing columns (G:L) that are no longer part of the sheet
$ws.Range("G1:L1").Clear()

$ws.Range("A1:F1").Select()
